$wb = $excel.ActiveWorkbook

# --- Step 1: turn the single "Sheet1" (OPER data) into two sheets:
#     MENU (new, placed first)  and  OPER (original data, placed second) ---
$orig = $wb.Worksheets.Item(1)
$orig.Copy($orig)

$menu = $wb.Worksheets.Item(1)
$oper = $wb.Worksheets.Item(2)

$menu.Name = "MENU"
$oper.Name = "OPER"

# --- Step 2: wipe the copied content on MENU, we'll repopulate it ---
$menu.Cells.Clear()

# --- Step 3: populate the MENU sheet layout ---
# Row 1: column index header 0..15 in B1:Q1 (bold + centered, like the
# original "B1 s=2" header style)
$menu.Range("B1").Value = 0
$menu.Range("C1").Value = 1
$menu.Range("D1").Value = 2
$menu.Range("E1").Value = 3
$menu.Range("F1").Value = 4
$menu.Range("G1").Value = 5
$menu.Range("H1").Value = 6
$menu.Range("I1").Value = 7
$menu.Range("J1").Value = 8
$menu.Range("K1").Value = 9
$menu.Range("L1").Value = 10
$menu.Range("M1").Value = 11
$menu.Range("N1").Value = 12
$menu.Range("O1").Value = 13
$menu.Range("P1").Value = 14
$menu.Range("Q1").Value = 15
$menu.Range("B1:Q1").Font.Bold = $true
$menu.Range("B1:Q1").HorizontalAlignment = -4108
$menu.Range("B1:Q1").VerticalAlignment = -4108

# Row 2: "TIMESLOT GPS "
$menu.Range("A2").Value = 0
$menu.Range("A2").Font.Bold = $true
$menu.Range("A2").HorizontalAlignment = -4108
$menu.Range("A2").VerticalAlignment = -4108
$menu.Range("B2").Value = "T"
$menu.Range("C2").Value = "I"
$menu.Range("D2").Value = "M"
$menu.Range("E2").Value = "E"
$menu.Range("F2").Value = "S"
$menu.Range("G2").Value = "L"
$menu.Range("H2").Value = "O"
$menu.Range("I2").Value = "T"
$menu.Range("J2").Value = " "
$menu.Range("K2").Value = "G"
$menu.Range("L2").Value = "P"
$menu.Range("M2").Value = "S"
$menu.Range("N2").Value = " "

# Row 3: "NEXT SEL UP"
$menu.Range("A3").Value = 1
$menu.Range("A3").Font.Bold = $true
$menu.Range("A3").HorizontalAlignment = -4108
$menu.Range("A3").VerticalAlignment = -4108
$menu.Range("B3").Value = "N"
$menu.Range("C3").Value = "E"
$menu.Range("D3").Value = "X"
$menu.Range("E3").Value = "T"
$menu.Range("G3").Value = "S"
$menu.Range("H3").Value = "E"
$menu.Range("I3").Value = "L"
$menu.Range("J3").Value = " "
$menu.Range("K3").Value = "U"
$menu.Range("L3").Value = "P"

# Row 7: slot-id style readout
$menu.Range("A7").Value = 0
$menu.Range("A7").Font.Bold = $true
$menu.Range("A7").HorizontalAlignment = -4108
$menu.Range("A7").VerticalAlignment = -4108
$menu.Range("B7").Value = 4
$menu.Range("C7").Value = "U"
$menu.Range("D7").Value = 1
$menu.Range("E7").Value = "U"
$menu.Range("F7").Value = "N"

# Row 8: same "NEXT SEL UP" labels as row 3
$menu.Range("A8").Value = 1
$menu.Range("A8").Font.Bold = $true
$menu.Range("A8").HorizontalAlignment = -4108
$menu.Range("A8").VerticalAlignment = -4108
$menu.Range("B8").Value = "N"
$menu.Range("C8").Value = "E"
$menu.Range("D8").Value = "X"
$menu.Range("E8").Value = "T"
$menu.Range("G8").Value = "S"
$menu.Range("H8").Value = "E"
$menu.Range("I8").Value = "L"
$menu.Range("J8").Value = " "
$menu.Range("K8").Value = "U"
$menu.Range("L8").Value = "P"

# Row 10
$menu.Range("A10").Value = 1
$menu.Range("A10").Font.Bold = $true
$menu.Range("A10").HorizontalAlignment = -4108
$menu.Range("A10").VerticalAlignment = -4108

# Row 13
$menu.Range("A13").Value = 1

# Row 15
$menu.Range("A15").Value = 1

# --- Step 4: selections / active sheet ---
$oper.Range("N3").Select()
$menu.Range("C11").Select()
$menu.Activate()
